$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A281").Value = 340
$ws.Range("B281").Value = "L1_DoubleJet40er2p5"
$ws.Range("D281").Value = 8400
$ws.Range("E281").Value = 8400
$ws.Range("F281").Value = 8400
$ws.Range("G281").Value = 8400
$ws.Range("H281").Value = 8400
$ws.Range("I281").Value = 8400
$ws.Range("J281").Value = 8400
$ws.Range("K281").Value = 8400

$ws.Range("A282").Value = 341
$ws.Range("B282").Value = "L1_DoubleJet100er2p5"
$ws.Range("D282").Value = 170
$ws.Range("E282").Value = 170
$ws.Range("F282").Value = 170
$ws.Range("G282").Value = 170
$ws.Range("H282").Value = 170
$ws.Range("I282").Value = 170
$ws.Range("J282").Value = 170
$ws.Range("K282").Value = 170

$ws.Range("A283").Value = 342
$ws.Range("B283").Value = "L1_DoubleJet120er2p5"
$ws.Range("D283").Value = 70
$ws.Range("E283").Value = 70
$ws.Range("F283").Value = 70
$ws.Range("G283").Value = 70
$ws.Range("H283").Value = 70
$ws.Range("I283").Value = 70
$ws.Range("J283").Value = 70
$ws.Range("K283").Value = 70

$ws.Range("A284").Value = 343
$ws.Range("B284").Value = "L1_DoubleJet150er2p5"
$ws.Range("D284").Value = 1
$ws.Range("E284").Value = 1
$ws.Range("F284").Value = 1
$ws.Range("G284").Value = 1
$ws.Range("H284").Value = 1
$ws.Range("I284").Value = 1
$ws.Range("J284").Value = 1
$ws.Range("K284").Value = 1

$ws.Range("A285").Value = 345
$ws.Range("B285").Value = "L1_DoubleJet100er2p3_dEta_Max1p6"
$ws.Range("D285").Value = 0
$ws.Range("E285").Value = 0
$ws.Range("F285").Value = 0
$ws.Range("G285").Value = 0
$ws.Range("H285").Value = 0
$ws.Range("I285").Value = 0
$ws.Range("J285").Value = 0
$ws.Range("K285").Value = 0

$ws.Range("A286").Value = 346
$ws.Range("B286").Value = "L1_DoubleJet112er2p3_dEta_Max1p6"
$ws.Range("D286").Value = 1
$ws.Range("E286").Value = 1
$ws.Range("F286").Value = 1
$ws.Range("G286").Value = 1
$ws.Range("H286").Value = 1
$ws.Range("I286").Value = 1
$ws.Range("J286").Value = 1
$ws.Range("K286").Value = 1

$ws.Range("A287").Value = 348
$ws.Range("B287").Value = "L1_DoubleJet30er2p5_Mass_Min150_dEta_Max1p5"
$ws.Range("D287").Value = 0
$ws.Range("E287").Value = 0
$ws.Range("F287").Value = 0
$ws.Range("G287").Value = 0
$ws.Range("H287").Value = 0
$ws.Range("I287").Value = 0
$ws.Range("J287").Value = 0
$ws.Range("K287").Value = 0

$ws.Range("A288").Value = 349
$ws.Range("B288").Value = "L1_DoubleJet30er2p5_Mass_Min200_dEta_Max1p5"
$ws.Range("D288").Value = 0
$ws.Range("E288").Value = 0
$ws.Range("F288").Value = 0
$ws.Range("G288").Value = 0
$ws.Range("H288").Value = 0
$ws.Range("I288").Value = 0
$ws.Range("J288").Value = 0
$ws.Range("K288").Value = 0

$ws.Range("A289").Value = 350
$ws.Range("B289").Value = "L1_DoubleJet30er2p5_Mass_Min250_dEta_Max1p5"
$ws.Range("D289").Value = 0
$ws.Range("E289").Value = 0
$ws.Range("F289").Value = 0
$ws.Range("G289").Value = 0
$ws.Range("H289").Value = 0
$ws.Range("I289").Value = 0
$ws.Range("J289").Value = 0
$ws.Range("K289").Value = 0

$ws.Range("A290").Value = 351
$ws.Range("B290").Value = "L1_DoubleJet30er2p5_Mass_Min300_dEta_Max1p5"
$ws.Range("D290").Value = 1
$ws.Range("E290").Value = 1
$ws.Range("F290").Value = 1
$ws.Range("G290").Value = 1
$ws.Range("H290").Value = 1
$ws.Range("I290").Value = 1
$ws.Range("J290").Value = 1
$ws.Range("K290").Value = 1

$ws.Range("A291").Value = 352
$ws.Range("B291").Value = "L1_DoubleJet30er2p5_Mass_Min330_dEta_Max1p5"
$ws.Range("D291").Value = 1
$ws.Range("E291").Value = 1
$ws.Range("F291").Value = 1
$ws.Range("G291").Value = 1
$ws.Range("H291").Value = 1
$ws.Range("I291").Value = 1
$ws.Range("J291").Value = 1
$ws.Range("K291").Value = 1

$ws.Range("A292").Value = 353
$ws.Range("B292").Value = "L1_DoubleJet30er2p5_Mass_Min360_dEta_Max1p5"
$ws.Range("D292").Value = 1
$ws.Range("E292").Value = 1
$ws.Range("F292").Value = 1
$ws.Range("G292").Value = 1
$ws.Range("H292").Value = 1
$ws.Range("I292").Value = 1
$ws.Range("J292").Value = 1
$ws.Range("K292").Value = 1

$ws.Range("A293").Value = 355
$ws.Range("B293").Value = "L1_DoubleJet_90_30_DoubleJet30_Mass_Min620"
$ws.Range("D293").Value = 0
$ws.Range("E293").Value = 0
$ws.Range("F293").Value = 0
$ws.Range("G293").Value = 0
$ws.Range("H293").Value = 0
$ws.Range("I293").Value = 0
$ws.Range("J293").Value = 0
$ws.Range("K293").Value = 0

$ws.Range("A294").Value = 356
$ws.Range("B294").Value = "L1_DoubleJet_100_30_DoubleJet30_Mass_Min620"
$ws.Range("D294").Value = 0
$ws.Range("E294").Value = 0
$ws.Range("F294").Value = 0
$ws.Range("G294").Value = 0
$ws.Range("H294").Value = 0
$ws.Range("I294").Value = 0
$ws.Range("J294").Value = 0
$ws.Range("K294").Value = 0

$ws.Range("A295").Value = 357
$ws.Range("B295").Value = "L1_DoubleJet_110_35_DoubleJet35_Mass_Min620"
$ws.Range("D295").Value = 0
$ws.Range("E295").Value = 1
$ws.Range("F295").Value = 1
$ws.Range("G295").Value = 1
$ws.Range("H295").Value = 1
$ws.Range("I295").Value = 1
$ws.Range("J295").Value = 1
$ws.Range("K295").Value = 1

$ws.Range("A296").Value = 358
$ws.Range("B296").Value = "L1_DoubleJet_115_40_DoubleJet40_Mass_Min620"
$ws.Range("D296").Value = 1
$ws.Range("E296").Value = 1
$ws.Range("F296").Value = 1
$ws.Range("G296").Value = 1
$ws.Range("H296").Value = 1
$ws.Range("I296").Value = 1
$ws.Range("J296").Value = 1
$ws.Range("K296").Value = 1

$ws.Range("A297").Value = 359
$ws.Range("B297").Value = "L1_DoubleJet_120_45_DoubleJet45_Mass_Min620"
$ws.Range("D297").Value = 1
$ws.Range("E297").Value = 1
$ws.Range("F297").Value = 1
$ws.Range("G297").Value = 1
$ws.Range("H297").Value = 1
$ws.Range("I297").Value = 1
$ws.Range("J297").Value = 1
$ws.Range("K297").Value = 1

$ws.Range("A298").Value = 360
$ws.Range("B298").Value = "L1_DoubleJet_115_40_DoubleJet40_Mass_Min620_Jet60TT28"
$ws.Range("D298").Value = 1
$ws.Range("E298").Value = 1
$ws.Range("F298").Value = 1
$ws.Range("G298").Value = 1
$ws.Range("H298").Value = 1
$ws.Range("I298").Value = 1
$ws.Range("J298").Value = 1
$ws.Range("K298").Value = 1

$ws.Range("A299").Value = 361
$ws.Range("B299").Value = "L1_DoubleJet_120_45_DoubleJet45_Mass_Min620_Jet60TT28"
$ws.Range("D299").Value = 1
$ws.Range("E299").Value = 1
$ws.Range("F299").Value = 1
$ws.Range("G299").Value = 1
$ws.Range("H299").Value = 1
$ws.Range("I299").Value = 1
$ws.Range("J299").Value = 1
$ws.Range("K299").Value = 1

$ws.Range("A300").Value = 362
$ws.Range("B300").Value = "L1_DoubleJet35_Mass_Min450_IsoTau45er2p1_RmOvlp_dR0p5"
$ws.Range("D300").Value = 0
$ws.Range("E300").Value = 0
$ws.Range("F300").Value = 0
$ws.Range("G300").Value = 0
$ws.Range("H300").Value = 0
$ws.Range("I300").Value = 0
$ws.Range("J300").Value = 0
$ws.Range("K300").Value = 0

$ws.Range("A301").Value = 363
$ws.Range("B301").Value = "L1_DoubleJet35_Mass_Min450_IsoTau45_RmOvlp"
$ws.Range("D301").Value = 0
$ws.Range("E301").Value = 0
$ws.Range("F301").Value = 0
$ws.Range("G301").Value = 0
$ws.Range("H301").Value = 0
$ws.Range("I301").Value = 0
$ws.Range("J301").Value = 0
$ws.Range("K301").Value = 0

$ws.Range("A302").Value = 364
$ws.Range("B302").Value = "L1_DoubleJet_80_30_Mass_Min420_IsoTau40_RmOvlp"
$ws.Range("D302").Value = 0
$ws.Range("E302").Value = 0
$ws.Range("F302").Value = 0
$ws.Range("G302").Value = 0
$ws.Range("H302").Value = 0
$ws.Range("I302").Value = 0
$ws.Range("J302").Value = 0
$ws.Range("K302").Value = 0

$ws.Range("A303").Value = 365
$ws.Range("B303").Value = "L1_DoubleJet_80_30_Mass_Min420_Mu8"
$ws.Range("D303").Value = 0
$ws.Range("E303").Value = 0
$ws.Range("F303").Value = 0
$ws.Range("G303").Value = 0
$ws.Range("H303").Value = 0
$ws.Range("I303").Value = 0
$ws.Range("J303").Value = 0
$ws.Range("K303").Value = 0

$ws.Range("A304").Value = 366
$ws.Range("B304").Value = "L1_DoubleJet_80_30_Mass_Min420_DoubleMu0_SQ"
$ws.Range("D304").Value = 0
$ws.Range("E304").Value = 0
$ws.Range("F304").Value = 0
$ws.Range("G304").Value = 0
$ws.Range("H304").Value = 0
$ws.Range("I304").Value = 0
$ws.Range("J304").Value = 0
$ws.Range("K304").Value = 0

$ws.Range("A305").Value = 372
$ws.Range("B305").Value = "L1_TripleJet_95_75_65_DoubleJet_75_65_er2p5"
$ws.Range("D305").Value = 0
$ws.Range("E305").Value = 1
$ws.Range("F305").Value = 1
$ws.Range("G305").Value = 1
$ws.Range("H305").Value = 1
$ws.Range("I305").Value = 1
$ws.Range("J305").Value = 1
$ws.Range("K305").Value = 1

$ws.Range("A306").Value = 373
$ws.Range("B306").Value = "L1_TripleJet_100_80_70_DoubleJet_80_70_er2p5"
$ws.Range("D306").Value = 1
$ws.Range("E306").Value = 1
$ws.Range("F306").Value = 1
$ws.Range("G306").Value = 1
$ws.Range("H306").Value = 1
$ws.Range("I306").Value = 1
$ws.Range("J306").Value = 1
$ws.Range("K306").Value = 1

$ws.Range("A307").Value = 374
$ws.Range("B307").Value = "L1_TripleJet_105_85_75_DoubleJet_85_75_er2p5"
$ws.Range("D307").Value = 1
$ws.Range("E307").Value = 1
$ws.Range("F307").Value = 1
$ws.Range("G307").Value = 1
$ws.Range("H307").Value = 1
$ws.Range("I307").Value = 1
$ws.Range("J307").Value = 1
$ws.Range("K307").Value = 1

$ws.Range("A308").Value = 376
$ws.Range("B308").Value = "L1_QuadJet_95_75_65_20_DoubleJet_75_65_er2p5_Jet20_FWD3p0"
$ws.Range("D308").Value = 1
$ws.Range("E308").Value = 1
$ws.Range("F308").Value = 1
$ws.Range("G308").Value = 1
$ws.Range("H308").Value = 1
$ws.Range("I308").Value = 1
$ws.Range("J308").Value = 1
$ws.Range("K308").Value = 1

$ws.Range("A309").Value = 382
$ws.Range("B309").Value = "L1_QuadJet60er2p5"
$ws.Range("D309").Value = 0
$ws.Range("E309").Value = 0
$ws.Range("F309").Value = 0
$ws.Range("G309").Value = 0
$ws.Range("H309").Value = 0
$ws.Range("I309").Value = 0
$ws.Range("J309").Value = 0
$ws.Range("K309").Value = 0

$ws.Range("A310").Value = 383
$ws.Range("B310").Value = "L1_HTT120_SingleLLPJet40"
$ws.Range("D310").Value = 0
$ws.Range("E310").Value = 0
$ws.Range("F310").Value = 0
$ws.Range("G310").Value = 0
$ws.Range("H310").Value = 0
$ws.Range("I310").Value = 0
$ws.Range("J310").Value = 0
$ws.Range("K310").Value = 0

$ws.Range("A311").Value = 384
$ws.Range("B311").Value = "L1_HTT160_SingleLLPJet50"
$ws.Range("D311").Value = 0
$ws.Range("E311").Value = 0
$ws.Range("F311").Value = 0
$ws.Range("G311").Value = 0
$ws.Range("H311").Value = 0
$ws.Range("I311").Value = 0
$ws.Range("J311").Value = 0
$ws.Range("K311").Value = 0

$ws.Range("A312").Value = 385
$ws.Range("B312").Value = "L1_HTT200_SingleLLPJet60"
$ws.Range("D312").Value = 0
$ws.Range("E312").Value = 0
$ws.Range("F312").Value = 0
$ws.Range("G312").Value = 0
$ws.Range("H312").Value = 0
$ws.Range("I312").Value = 0
$ws.Range("J312").Value = 0
$ws.Range("K312").Value = 0

$ws.Range("A313").Value = 386
$ws.Range("B313").Value = "L1_HTT240_SingleLLPJet70"
$ws.Range("D313").Value = 0
$ws.Range("E313").Value = 0
$ws.Range("F313").Value = 0
$ws.Range("G313").Value = 0
$ws.Range("H313").Value = 0
$ws.Range("I313").Value = 0
$ws.Range("J313").Value = 0
$ws.Range("K313").Value = 0

$ws.Range("A314").Value = 387
$ws.Range("B314").Value = "L1_DoubleLLPJet40"
$ws.Range("D314").Value = 0
$ws.Range("E314").Value = 0
$ws.Range("F314").Value = 0
$ws.Range("G314").Value = 0
$ws.Range("H314").Value = 0
$ws.Range("I314").Value = 0
$ws.Range("J314").Value = 0
$ws.Range("K314").Value = 0

$ws.Range("A315").Value = 388
$ws.Range("B315").Value = "L1_HTT280er_QuadJet_70_55_40_35_er2p5"
$ws.Range("D315").Value = 0
$ws.Range("E315").Value = 0
$ws.Range("F315").Value = 0
$ws.Range("G315").Value = 0
$ws.Range("H315").Value = 0
$ws.Range("I315").Value = 0
$ws.Range("J315").Value = 0
$ws.Range("K315").Value = 0

$ws.Range("A316").Value = 389
$ws.Range("B316").Value = "L1_HTT320er_QuadJet_70_55_40_40_er2p5"
$ws.Range("D316").Value = 0
$ws.Range("E316").Value = 1
$ws.Range("F316").Value = 1
$ws.Range("G316").Value = 1
$ws.Range("H316").Value = 1
$ws.Range("I316").Value = 1
$ws.Range("J316").Value = 1
$ws.Range("K316").Value = 1

$ws.Range("A317").Value = 390
$ws.Range("B317").Value = "L1_HTT320er_QuadJet_80_60_er2p1_45_40_er2p3"
$ws.Range("D317").Value = 1
$ws.Range("E317").Value = 1
$ws.Range("F317").Value = 1
$ws.Range("G317").Value = 1
$ws.Range("H317").Value = 1
$ws.Range("I317").Value = 1
$ws.Range("J317").Value = 1
$ws.Range("K317").Value = 1

$ws.Range("A318").Value = 391
$ws.Range("B318").Value = "L1_HTT320er_QuadJet_80_60_er2p1_50_45_er2p3"
$ws.Range("D318").Value = 1
$ws.Range("E318").Value = 1
$ws.Range("F318").Value = 1
$ws.Range("G318").Value = 1
$ws.Range("H318").Value = 1
$ws.Range("I318").Value = 1
$ws.Range("J318").Value = 1
$ws.Range("K318").Value = 1
